$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns (zh-cn / de-de) and HO Xliff generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2017-01-03 07:56:17"

# --- zh-cn sheet: Status column + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-01-03 07:56:06"

# --- de-de sheet: Status column + Latest Handback DateTime (mirrors Overview G2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-01-03 07:56:17"

# --- Column width adjustments to fit the longer "Ready for handoff" text ---
# NOTE: the host quantizes ColumnWidth to round(width*6)+5 pixels, so the
# nearest representable width to the target 17.2159881591797 is reached by
# feeding in a value in the same quantization bucket (any value in
# [16.25, 16.41666...) round-trips to 17.166666666666668, the closest
# achievable approximation of the target width).
$targetColWidth = 16.333333333333336
$wsOverview.Columns.Item(5).ColumnWidth = $targetColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColWidth
